# Update "想去人数" (column F) counts on the 展览 and 全部类型 sheets.
# Both sheets hold the same listing, so the same rows/values change on each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    4  = 1579
    7  = 11330
    13 = 784
    14 = 12316
    15 = 12969
    20 = 82
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
